{"js": "// Replace the date and each \"a\u00f7b=c, r\" table-cell value with its new value.\n// Every old text value in this document is unique, so a plain search +\n// replace (matchCase, whole-string) for each pair is safe and unambiguous.\nconst replacements = [\n  [\"2025-09-01 Monday\", \"2025-09-02 Tuesday\"],\n  [\"556\u00f77=79, 3\", \"702\u00f76=117, 0\"],\n  [\"182\u00f77=26, 0\", \"332\u00f72=166, 0\"],\n  [\"504\u00f79=56, 0\", \"924\u00f74=231, 0\"],\n  [\"775\u00f79=86, 1\", \"957\u00f78=119, 5\"],\n  [\"511\u00f73=170, 1\", \"230\u00f76=38, 2\"],\n  [\"749\u00f75=149, 4\", \"529\u00f72=264, 1\"],\n  [\"988\u00f75=197, 3\", \"685\u00f74=171, 1\"],\n  [\"128\u00f79=14, 2\", \"352\u00f78=44, 0\"],\n  [\"390\u00f77=55, 5\", \"370\u00f74=92, 2\"],\n  [\"131\u00f79=14, 5\", \"238\u00f73=79, 1\"],\n  [\"459\u00f77=65, 4\", \"746\u00f72=373, 0\"],\n  [\"974\u00f76=162, 2\", \"257\u00f78=32, 1\"],\n  [\"621\u00f72=310, 1\", \"107\u00f76=17, 5\"],\n  [\"275\u00f78=34, 3\", \"365\u00f76=60, 5\"],\n  [\"167\u00f78=20, 7\", \"318\u00f77=45, 3\"],\n  [\"798\u00f77=114, 0\", \"985\u00f72=492, 1\"],\n  [\"922\u00f78=115, 2\", \"573\u00f74=143, 1\"],\n  [\"886\u00f72=443, 0\", \"464\u00f75=92, 4\"],\n  [\"894\u00f79=99, 3\", \"521\u00f78=65, 1\"],\n  [\"436\u00f72=218, 0\", \"634\u00f78=79, 2\"],\n  [\"374\u00f73=124, 2\", \"228\u00f76=38, 0\"],\n  [\"931\u00f73=310, 1\", \"354\u00f72=177, 0\"],\n  [\"541\u00f73=180, 1\", \"573\u00f78=71, 5\"],\n  [\"158\u00f78=19, 6\", \"399\u00f72=199, 1\"],\n  [\"203\u00f78=25, 3\", \"782\u00f78=97, 6\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"a\u00f7b=c, r\" table-cell value with its new value.\n# Every old text value in this document is unique, so a plain Find/Replace\n# (MatchCase, MatchWholeWord) for each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-01 Monday\", \"2025-09-02 Tuesday\"),\n    @(\"556\u00f77=79, 3\", \"702\u00f76=117, 0\"),\n    @(\"182\u00f77=26, 0\", \"332\u00f72=166, 0\"),\n    @(\"504\u00f79=56, 0\", \"924\u00f74=231, 0\"),\n    @(\"775\u00f79=86, 1\", \"957\u00f78=119, 5\"),\n    @(\"511\u00f73=170, 1\", \"230\u00f76=38, 2\"),\n    @(\"749\u00f75=149, 4\", \"529\u00f72=264, 1\"),\n    @(\"988\u00f75=197, 3\", \"685\u00f74=171, 1\"),\n    @(\"128\u00f79=14, 2\", \"352\u00f78=44, 0\"),\n    @(\"390\u00f77=55, 5\", \"370\u00f74=92, 2\"),\n    @(\"131\u00f79=14, 5\", \"238\u00f73=79, 1\"),\n    @(\"459\u00f77=65, 4\", \"746\u00f72=373, 0\"),\n    @(\"974\u00f76=162, 2\", \"257\u00f78=32, 1\"),\n    @(\"621\u00f72=310, 1\", \"107\u00f76=17, 5\"),\n    @(\"275\u00f78=34, 3\", \"365\u00f76=60, 5\"),\n    @(\"167\u00f78=20, 7\", \"318\u00f77=45, 3\"),\n    @(\"798\u00f77=114, 0\", \"985\u00f72=492, 1\"),\n    @(\"922\u00f78=115, 2\", \"573\u00f74=143, 1\"),\n    @(\"886\u00f72=443, 0\", \"464\u00f75=92, 4\"),\n    @(\"894\u00f79=99, 3\", \"521\u00f78=65, 1\"),\n    @(\"436\u00f72=218, 0\", \"634\u00f78=79, 2\"),\n    @(\"374\u00f73=124, 2\", \"228\u00f76=38, 0\"),\n    @(\"931\u00f73=310, 1\", \"354\u00f72=177, 0\"),\n    @(\"541\u00f73=180, 1\", \"573\u00f78=71, 5\"),\n    @(\"158\u00f78=19, 6\", \"399\u00f72=199, 1\"),\n    @(\"203\u00f78=25, 3\", \"782\u00f78=97, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)  # 2 = wdReplaceAll\n}\n"}
